$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G6 already carries the "no highlight" style used elsewhere on this sheet.
# Copy its format onto C6:H6 (replacing the green highlight C6:F6 currently
# have) and fill every one of those six cells in with a score of 5.
$ws.Range("G6").Copy()
$ws.Range("C6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C6:H6").Value2 = 5

# Move the active selection to I6.
$ws.Range("I6").Select()
